# Adds changes in relation to P drive validations.
# On the "Constants" sheet:
#  - Row 25 (PathTempToPDF) value changes from the local P:\ drive path to the
#    UNC equivalent \\10.250.52.158\Depts\...
#  - A new row is inserted right after it (new row 26) defining
#    PathDifferentTempToPDF with a second UNC fallback path
#    (\\somproddfs1.prod.sovos.org\depts\...), mirroring the existing
#    PathPDriveFolder / PathDifferentPDriveFolder pair above it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Constants")

# Insert a new row below the "PathTempToPDF" row (row 25) so the new setting
# lines up right under it, shifting everything below down by one.
$ws.Rows.Item(26).Insert()
$ws.Rows.Item(26).RowHeight = 14.25

# Update the PathTempToPDF value to use the UNC path instead of the mapped
# P:\ drive path.
$ws.Range("B25").Value = "\\10.250.52.158\Depts\TaxReturnOutSourcing\Preparer\UIPathPublish\IR Bot Temp Files\pdfTemp <USERNAME> <CUSTOMER>.pdf"
$ws.Range("B25").Interior.ColorIndex = -4142

# Populate the new row with the alternate UNC path setting.
$ws.Range("A26").Value = "PathDifferentTempToPDF"
$ws.Range("B26").Value = "\\somproddfs1.prod.sovos.org\depts\TaxReturnOutSourcing\Preparer\UIPathPublish\IR Bot Temp Files\pdfTemp <USERNAME> <CUSTOMER>.pdf"
$ws.Range("B26").Interior.ColorIndex = -4142

# Leave the sheet focused on the area that was just edited.
$ws.Activate() | Out-Null
$ws.Range("A26").Select() | Out-Null
